$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 4-9 (previously blank) with station/date/time data,
# mirroring the pattern already present in rows 2-3.
$ws.Range("B4").Value = "HINDUJA"
$ws.Range("C4").Value = 46024
$ws.Range("D4").Value = 0.0729166666666667
$ws.Range("E4").Value = 46024
$ws.Range("F4").Value = 0.125

$ws.Range("B5").Value = "HINDUJA"
$ws.Range("C5").Value = 46024
$ws.Range("D5").Value = 0.916666666666667
$ws.Range("E5").Value = 46024
$ws.Range("F5").Value = 1

$ws.Range("B6").Value = "HINDUJA"
$ws.Range("C6").Value = 46025
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 46025
$ws.Range("F6").Value = 0.138888888888889

$ws.Range("B7").Value = "HINDUJA"
$ws.Range("C7").Value = 46025
$ws.Range("D7").Value = 0.4375
$ws.Range("E7").Value = 46025
$ws.Range("F7").Value = 0.572916666666667

$ws.Range("B8").Value = "HINDUJA"
$ws.Range("C8").Value = 46025
$ws.Range("D8").Value = 0.84375
$ws.Range("E8").Value = 46025
$ws.Range("F8").Value = 0.861111111111111

$ws.Range("B9").Value = "HINDUJA"
$ws.Range("C9").Value = 46025
$ws.Range("D9").Value = 0.920138888888889
$ws.Range("E9").Value = 46025
$ws.Range("F9").Value = 1

# New column G mirrors column F's formatting (time, boxed) for the rows
# that now carry data, but is left without a value for now.
$ws.Range("F2:F9").Copy()
$ws.Range("G2:G9").PasteSpecial(-4122)

# The S.No column (A) is cleared for all of the still-unused rows below
# row 9 (previously pre-numbered 9..99) so the sheet only numbers
# entered records.
$ws.Range("A10:A100").ClearContents()

# Move/collapse the active selection onto the next entry row.
$ws.Range("B16").Select() | Out-Null
